$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the policy number (NumPoliza) in row 2 from 04104009144 to 04104013194
# Keep it stored as text (it has a leading-zero numeric-looking value), matching
# the existing quote-prefixed style already applied to this cell.
$ws.Range("E2").Value = "'04104013194"

# Update selection to mirror the recorded cursor position after the edit
$ws.Range("F16").Select()
